# EMEP_NFR14_scaling_mapping.xlsx - add new country scaling rows to the
# "year" sheet (kgz, arm, blr, mlt), matching the new India (Venkataraman)
# scaling inventory commit.

$wb = $excel.ActiveWorkbook
$wsMap   = $wb.Worksheets.Item("map")
$wsYear  = $wb.Worksheets.Item("year")

# --- New data rows on the "year" sheet (rows 5-8) ---------------------
$wsYear.Range("A5").Value = "kgz"
$wsYear.Range("B5").Value = "all"
$wsYear.Range("C5").Value = "NA"
$wsYear.Range("D5").Value = "NA"
$wsYear.Range("E5").Value = "NA"
$wsYear.Range("F5").Value = 2016
$wsYear.Range("G5").Value = 2017
$wsYear.Range("H5").Value = "avoid missing data years"

$wsYear.Range("A6").Value = "arm"
$wsYear.Range("B6").Value = "all"
$wsYear.Range("C6").Value = "NA"
$wsYear.Range("D6").Value = "NA"
$wsYear.Range("E6").Value = "NA"
$wsYear.Range("F6").Value = 2016
$wsYear.Range("G6").Value = 2017

$wsYear.Range("A7").Value = "blr"
$wsYear.Range("B7").Value = "all"
$wsYear.Range("C7").Value = "NA"
$wsYear.Range("D7").Value = "NA"
$wsYear.Range("E7").Value = "NA"
$wsYear.Range("F7").Value = 2014
$wsYear.Range("G7").Value = 2017

$wsYear.Range("A8").Value = "mlt"
$wsYear.Range("B8").Value = "all"
$wsYear.Range("C8").Value = "NA"
$wsYear.Range("D8").Value = "NA"
$wsYear.Range("E8").Value = "NA"
$wsYear.Range("F8").Value = 2000
$wsYear.Range("G8").Value = 2017

# --- View-state changes: "map" loses the active tab/selection, -------
# --- "year" becomes the active sheet with a fresh selection. ---------
$wsMap.Range("A31").Select()
$wsYear.Activate()
$wsYear.Range("A9").Select()
